# Fusion-Tracker.xlsx — "added task 11 and 12"
#
# Summary sheet: append Task 11 / Task 12 rows (21-22).
# Daily Updates sheet: append the three daily-log rows (14-16) documenting
# the work done for Task 11 and Task 12.

$wb = $excel.ActiveWorkbook

$dateFmt = "[$-409]d/mmm/yyyy;@"

# ----- Summary sheet -----
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("A21").Value = "Task 11"
$summary.Range("B21").Value = 44379
$summary.Range("B21").NumberFormat = $dateFmt
$summary.Range("C21").Value = "Bi - Purchase Order Detailed Report"

$summary.Range("A22").Value = "Task 12"
$summary.Range("B22").Value = 44383
$summary.Range("B22").NumberFormat = $dateFmt
$summary.Range("C22").Value = "SQL - Joins, Set Operators and Functions"
$summary.Range("D22").Value = 44383
$summary.Range("D22").NumberFormat = $dateFmt
$summary.Range("E22").Value = "Y"

$summary.Range("E22").Select()

# ----- Daily Updates sheet -----
$daily = $wb.Worksheets.Item("Daily Updates")

$daily.Range("A14").Value = 44379
$daily.Range("A14").NumberFormat = $dateFmt
$daily.Range("B14").Value = "Task 11"
$daily.Range("C14").Value = "FSD analysis"
$daily.Range("D14").Value = "Table and Column identification"
$daily.Range("E14").Value = "Task 11 analysed"

$daily.Range("A15").Value = 44382
$daily.Range("A15").NumberFormat = $dateFmt
$daily.Range("B15").Value = "Task 11"
$daily.Range("C15").Value = "Creation of data model"
$daily.Range("D15").Value = "Generating EXCEL template"
$daily.Range("E15").Value = "Task 11 template generated"

$daily.Range("A16").Value = 44383
$daily.Range("A16").NumberFormat = $dateFmt
$daily.Range("B16").Value = "Task 12"
$daily.Range("C16").Value = "Properties - Data Model and Report"
$daily.Range("D16").Value = "SQL Documentation"
$daily.Range("E16").Value = "Task 12 completed"

$daily.Range("E16").Select()

# Re-activate Summary so the workbook reopens on the same tab as before.
$summary.Activate()
